$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = -7
$ws.Range("F3").Value = 3
$ws.Range("F8").Value = -4
$ws.Range("F14").Value = -2
$ws.Range("F16").Value = 1
$ws.Range("F17").Value = 0
$ws.Range("F18").Value = -1
$ws.Range("F22").Value = -1
